$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FINE-TUNING")

# --- Row 10: new fine-tuning model "dots-llmv6" results ---
$ws.Range("H10").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv6:B4Y49YIK:ckpt-step-70"
$ws.Range("I10").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv6:B4Y493eN:ckpt-step-84"
$ws.Range("J10").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv6:B4Y4AFah"
$ws.Range("K10").Value = 250663
$ws.Range("L10").Value = 7
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = 1.8
$ws.Range("O10").Value = 191679768

# Checkpoint names (H10:J10) are centered, no wrap
$ws.Range("H10:J10").HorizontalAlignment = -4108
$ws.Range("H10:J10").VerticalAlignment = -4108
$ws.Range("H10:J10").WrapText = $false

# Epochs / batch size / LR multiplier (L10:N10) match the header styling (center + wrap)
$ws.Range("L10:N10").HorizontalAlignment = -4108
$ws.Range("L10:N10").VerticalAlignment = -4108
$ws.Range("L10:N10").WrapText = $true

# --- Move the active selection to H10, as in the saved workbook ---
$ws.Activate() | Out-Null
$ws.Range("H10").Select() | Out-Null
